# Update column F (dSF) values per repulled data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -2
    5  = 1
    6  = 2
    8  = 4
    9  = 1
    10 = 1
    12 = 3
    13 = -1
    15 = 1
    16 = 1
    17 = 1
    20 = -1
    21 = 1
    22 = 3
    23 = 1
    24 = 2
    25 = -4
    26 = -1
    27 = -3
    28 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
